$wb = $excel.ActiveWorkbook

# --- Weekly Quantity sheet: remove the row for 2024-03-10 (value 45361.99999999999 / 420) ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(15).Delete()

# --- Monthly Trend sheet: correct requested quantity for the 2024-03 row from 1746 to 1326 ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(6, 2).Value = 1326
